$d = $word.ActiveDocument

# Update the date heading paragraph.
$d.Content.Find.Execute("2024-10-26 Saturday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-10-27 Sunday", 2) | Out-Null

# Update the division-problem table. Using direct cell addressing (rather than a
# global Find/Replace) avoids ambiguity since some new values coincide with other
# cells' old values (e.g. "21÷7=" is both an old value and a replacement result).
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "47÷7=" },
    @{ Row = 1;  Col = 2; Text = "66÷7=" },
    @{ Row = 1;  Col = 3; Text = "21÷7=" },
    @{ Row = 1;  Col = 4; Text = "68÷5=" },
    @{ Row = 1;  Col = 5; Text = "40÷3=" },

    @{ Row = 5;  Col = 1; Text = "67÷8=" },
    @{ Row = 5;  Col = 2; Text = "31÷4=" },
    @{ Row = 5;  Col = 3; Text = "33÷6=" },
    @{ Row = 5;  Col = 4; Text = "68÷9=" },
    @{ Row = 5;  Col = 5; Text = "64÷3=" },

    @{ Row = 9;  Col = 1; Text = "42÷5=" },
    @{ Row = 9;  Col = 2; Text = "64÷3=" },
    @{ Row = 9;  Col = 3; Text = "34÷9=" },
    @{ Row = 9;  Col = 4; Text = "81÷4=" },
    @{ Row = 9;  Col = 5; Text = "92÷9=" },

    @{ Row = 13; Col = 1; Text = "37÷2=" },
    @{ Row = 13; Col = 2; Text = "51÷7=" },
    @{ Row = 13; Col = 3; Text = "20÷6=" },
    @{ Row = 13; Col = 4; Text = "85÷2=" },
    @{ Row = 13; Col = 5; Text = "73÷9=" },

    @{ Row = 17; Col = 1; Text = "29÷7=" },
    @{ Row = 17; Col = 2; Text = "27÷9=" },
    @{ Row = 17; Col = 3; Text = "49÷2=" },
    @{ Row = 17; Col = 4; Text = "24÷2=" },
    @{ Row = 17; Col = 5; Text = "91÷4=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}

Write-Host "Done updating date and" $updates.Count "table cells."
